# "Arreglado diseño del proyecto"
#
# The slide is one big top-level group ("Grupo 8") that contains the whole
# AWS-IoT architecture diagram. Three of its nested sub-groups
# (cognitoLogs / getLastTemp / HeatSense) were repositioned, and the
# straight-connector arrows glued to them were re-routed to match.
#
# This COM host flattens nested p:grpSp containers: Shapes.Item(1) is the
# single outer group and its .GroupItems exposes every leaf shape
# (pic / sp / cxnSp) directly, regardless of how many grpSp wrappers they
# were nested inside. The intermediate "Grupo N" wrapper shapes themselves
# are not individually addressable objects here, so we recreate the effect
# of "moving a sub-group" by shifting each of its leaf children by the same
# delta, and we set the affected connectors' geometry explicitly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$top = $s.Shapes.Item(1)

function Get-ShapeById($container, [int]$id) {
    for ($i = 1; $i -le $container.GroupItems.Count; $i++) {
        $it = $container.GroupItems.Item($i)
        if ($it.Id -eq $id) { return $it }
    }
    return $null
}

# Cosmetic rename that went along with the fix (id itself is read-only).
$top.Name = "Grupo 2"

# --- Move the "cognitoLogs" sub-group (picture 75 + caption 76) -----------
$pic75 = Get-ShapeById $top 75
$pic75.Left = 228.6267716535433
$pic75.Top  = 295.6188188976378

$cap76 = Get-ShapeById $top 76
$cap76.Left = 205.84149606299212
$cap76.Top  = 344.4996062992126

# --- Move the "getLastTemp" sub-group (picture 88 + caption 89) -----------
$pic88 = Get-ShapeById $top 88
$pic88.Left = 382.0081102362205
$pic88.Top  = 36.8503937007874

$cap89 = Get-ShapeById $top 89
$cap89.Left = 358.92070866141734
$cap89.Top  = 91.56779527559056

# --- Move the "HeatSense" sub-group (picture 123 + caption 128) -----------
$pic123 = Get-ShapeById $top 123
$pic123.Left = 790.9448818897638
$pic123.Top  = 155.99818897637795

$cap128 = Get-ShapeById $top 128
$cap128.Left = 766.9918110236221
$cap128.Top  = 196.32125984251968

# --- Re-route the connector arrows that touch the moved shapes ------------

# "Conector recto de flecha 30" (saveLogs/Cognito user icon -> web)
$cxn31 = Get-ShapeById $top 31
$cxn31.Left = 616.5068503937008
$cxn31.Top = 350.1812598425197
$cxn31.Width = 51.57070866141732
$cxn31.Height = 68.38795275590552
$cxn31.HorizontalFlip = -1
$cxn31.VerticalFlip = 0

# "Conector recto de flecha 72" (saveLogs -> cognitoLogs)
$cxn73 = Get-ShapeById $top 73
$cxn73.Left = 539.7921259842519
$cxn73.Top = 222.9816535433071
$cxn73.Width = 22.341023622047246
$cxn73.Height = 22.976614173228345
$cxn73.HorizontalFlip = 0
$cxn73.VerticalFlip = 0

# "Conector recto de flecha 110" (cognitoLogs/lastTemp -> getLastTemp)
$cxn111 = Get-ShapeById $top 111
try {
    $pic49 = Get-ShapeById $top 49
    $cxn111.ConnectorFormat.BeginDisconnect()
    $cxn111.ConnectorFormat.BeginConnect($pic49, 3)
} catch {}
$cxn111.Left = 538.0727559055118
$cxn111.Top = 128.47181102362205
$cxn111.Width = 130.08448818897637
$cxn111.Height = 141.5459842519685
$cxn111.HorizontalFlip = 0
$cxn111.VerticalFlip = 0

# "Conector recto de flecha 123" (web -> HeatSense)
$cxn124 = Get-ShapeById $top 124
$cxn124.Left = 581.5041732283464
$cxn124.Top = 364.9617322834646
$cxn124.Width = 2.635748031496063
$cxn124.Height = 15.295511811023623
$cxn124.HorizontalFlip = -1
$cxn124.VerticalFlip = -1

# "Conector recto de flecha 130" (HeatSense -> getLastTemp)
$cxn131 = Get-ShapeById $top 131
$cxn131.Left = 597.9464566929133
$cxn131.Top = 297.1003937007874
$cxn131.Width = 46.430944881889765
$cxn131.Height = 22.253779527559054
$cxn131.HorizontalFlip = 0
$cxn131.VerticalFlip = -1
